# Express-test sheet: change value parsing (extra measured columns) and
# a small header/label fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("A1").Value = "Модуляция"
$ws.Range("B1").Value = "Аттен, ДБ"
$ws.Range("C1").Value = "С/Ш"
$ws.Range("D1").Value = "Отправлено, байт"
$ws.Range("E1").Value = "Принято, байт"
$ws.Range("F1").Value = "Потерято, байт"
$ws.Range("G1").Value = "Процент ошибок, %"

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = "QAM64 3/4"
$ws.Range("B2").Value = 20
# C2 keeps the "number stored as text" shape it had before the edit.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "26.350000381469727"
$ws.Range("D2").Value = 16646272
$ws.Range("E2").Value = 16347840
$ws.Range("F2").Value = 298432
$ws.Range("G2").Value = 1.7927857961229996

# The sheet's used range grew from A1:D2 to A1:G2; keep Excel's
# "ignore number stored as text" error-check region in sync with it.
$ws.Range("A1:G2").Errors(3).Ignore = $true
